$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"
$zhcnHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$dedeHandbackFile = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$zhcnHandbackDate = "2016-08-21 16:45:19"
$dedeHandbackDate = "2016-08-21 16:45:25"

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4b564cfe5d8fdf6ee07b2623916e877437318fc7/e2e/a.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4b564cfe5d8fdf6ee07b2623916e877437318fc7/e2e/b.md"

# Hyperlink-style blue color (BGR-encoded OLE_COLOR for RGB 6495ED)
$hyperlinkColor = 0xED9564

function Set-HandbackInfo($ws, $dedupDate) {
    # Update status on row 2 and row 3
    $ws.Range("C2").Value = $statusHandedBack
    $ws.Range("C3").Value = $statusHandedBack

    # Rebuild hyperlinks (A2, A3) plus add the two new ones (I2, I3)
    # so relationship ids come out in document order: A2, I2, A3, I3
    $ws.Range("A2").Hyperlinks.Delete()
    $ws.Range("A3").Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $urlB, [Type]::Missing, [Type]::Missing, "b.md")
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlA, [Type]::Missing, [Type]::Missing, "a.md")

    foreach ($addr in @("A2", "A3", "I2", "I3")) {
        $c = $ws.Range($addr)
        $c.Font.Underline = 2
        $c.Font.Color = $hyperlinkColor
    }
}

# zh-cn sheet
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Set-HandbackInfo $wsZhCn $zhcnHandbackDate
$wsZhCn.Range("J2").Value = $zhcnHandbackFile
$wsZhCn.Range("J3").Value = $zhcnHandbackFile
$wsZhCn.Range("K2").Value = $zhcnHandbackDate
$wsZhCn.Range("K3").Value = $zhcnHandbackDate
$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsZhCn.Columns.Item(10).ColumnWidth = 39.166666666666664

# de-de sheet
$wsDeDe = $wb.Worksheets.Item("de-de")
Set-HandbackInfo $wsDeDe $dedeHandbackDate
$wsDeDe.Range("J2").Value = $dedeHandbackFile
$wsDeDe.Range("J3").Value = $dedeHandbackFile
$wsDeDe.Range("K2").Value = $dedeHandbackDate
$wsDeDe.Range("K3").Value = $dedeHandbackDate
$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666664
$wsDeDe.Columns.Item(10).ColumnWidth = 39.166666666666664

# Overview sheet mirrors the status text in columns E (zh-cn) and F (de-de)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666664
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666664

Write-Output "Handback report generated"
